# Auto-generated Excel COM-interop script
# Applies numeric corrections to columns H:N across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 6059.4614
$arr[0,1] = 3110.5
$arr[0,2] = 6595.636
$arr[0,3] = 3110.5
$arr[0,4] = 6595.636
$ws.Range("H62:L62").Value = $arr
$ws.Range("M62").Value = -2486.5
$ws.Range("N62").Value = -7843.636

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 6059.4614
$arr[0,1] = 3110.5
$arr[0,2] = 6595.636
$arr[0,3] = 15552.5
$arr[0,4] = 32978.18
$ws.Range("H65:L65").Value = $arr
$ws.Range("M65").Value = -12432.5
$ws.Range("N65").Value = -39218.18

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 49078.25
$arr[0,1] = 41797
$arr[0,2] = 100047
$arr[0,3] = 125391
$arr[0,4] = 300141
$ws.Range("H82:L82").Value = $arr
$ws.Range("M82").Value = -124985
$ws.Range("N82").Value = -300953

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 49078.25
$arr[0,1] = 41797
$arr[0,2] = 100047
$arr[0,3] = 125391
$arr[0,4] = 300141
$ws.Range("H85:L85").Value = $arr
$ws.Range("M85").Value = -123987
$ws.Range("N85").Value = -302949

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 2140.4666
$arr[0,1] = 1342.3334
$arr[0,2] = 5333
$arr[0,3] = 4027.0002
$arr[0,4] = 15999
$ws.Range("H96:L96").Value = $arr
$ws.Range("M96").Value = -2654.0002
$ws.Range("N96").Value = -18745

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1070.1666
$arr[0,1] = 1043
$arr[0,2] = 1206
$arr[0,3] = 1043
$arr[0,4] = 1206
$ws.Range("H100:L100").Value = $arr
$ws.Range("M100").Value = -502
$ws.Range("N100").Value = -2288

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 827.6667
$arr[0,1] = 861.7857
$arr[0,2] = 350
$arr[0,3] = 861.7857
$arr[0,4] = 350
$ws.Range("H107:L107").Value = $arr
$ws.Range("M107").Value = 1058.2143
$ws.Range("N107").Value = -4190

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 126963.336
$arr[0,1] = 0
$arr[0,2] = 126963.336
$arr[0,3] = 0
$arr[0,4] = 126963.336
$ws.Range("H117:L117").Value = $arr
$ws.Range("N117").Value = -136141.336

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 274869.7
$arr[0,1] = 372811.84
$arr[0,2] = 10425.8
$arr[0,3] = 1118435.52
$arr[0,4] = 31277.4
$ws.Range("H137:L137").Value = $arr
$ws.Range("M137").Value = -1115885.52
$ws.Range("N137").Value = -36377.39999999999

$ws = $wb.Worksheets.Item("ARM")
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1408.3636
$arr[0,1] = 1419.4
$arr[0,2] = 1298
$arr[0,3] = 1419.4
$arr[0,4] = 1298
$ws.Range("H2:L2").Value = $arr
$ws.Range("M2").Value = -1306.4
$ws.Range("N2").Value = -1524

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 9293777
$arr[0,1] = 14743863
$arr[0,2] = 28630.2
$arr[0,3] = 14743863
$arr[0,4] = 28630.2
$ws.Range("H32:L32").Value = $arr
$ws.Range("M32").Value = -14743576
$ws.Range("N32").Value = -29204.2

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1488.1351
$arr[0,1] = 1342.3667
$arr[0,2] = 2112.8572
$arr[0,3] = 1342.3667
$arr[0,4] = 2112.8572
$ws.Range("H97:L97").Value = $arr
$ws.Range("M97").Value = -846.3667
$ws.Range("N97").Value = -3104.8572

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1408.3636
$arr[0,1] = 1419.4
$arr[0,2] = 1298
$arr[0,3] = 1419.4
$arr[0,4] = 1298
$ws.Range("H116:L116").Value = $arr
$ws.Range("M116").Value = 874.5999999999999
$ws.Range("N116").Value = -5886

$ws = $wb.Worksheets.Item("BSM")
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1408.3636
$arr[0,1] = 1419.4
$arr[0,2] = 1298
$arr[0,3] = 1419.4
$arr[0,4] = 1298
$ws.Range("H3:L3").Value = $arr
$ws.Range("M3").Value = -1305.4
$ws.Range("N3").Value = -1526

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 2319.8333
$arr[0,1] = 2201.5625
$arr[0,2] = 2556.375
$arr[0,3] = 2201.5625
$arr[0,4] = 2556.375
$ws.Range("H86:L86").Value = $arr
$ws.Range("M86").Value = -1078.5625
$ws.Range("N86").Value = -4802.375

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 2319.8333
$arr[0,1] = 2201.5625
$arr[0,2] = 2556.375
$arr[0,3] = 11007.8125
$arr[0,4] = 12781.875
$ws.Range("H89:L89").Value = $arr
$ws.Range("M89").Value = -5391.8125
$ws.Range("N89").Value = -24013.875

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 2466.7778
$arr[0,1] = 1621.909
$arr[0,2] = 3794.4285
$arr[0,3] = 1621.909
$arr[0,4] = 3794.4285
$ws.Range("H107:L107").Value = $arr
$ws.Range("M107").Value = 298.0909999999999
$ws.Range("N107").Value = -7634.4285

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 2130399.5
$arr[0,1] = 1667.2424
$arr[0,2] = 7148125.5
$arr[0,3] = 5001.7272
$arr[0,4] = 21444376.5
$ws.Range("H134:L134").Value = $arr
$ws.Range("M134").Value = -2466.7272
$ws.Range("N134").Value = -21449446.5

$ws = $wb.Worksheets.Item("CRP")
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 5965.3335
$arr[0,1] = 4881.2
$arr[0,2] = 7320.5
$arr[0,3] = 4881.2
$arr[0,4] = 7320.5
$ws.Range("H58:L58").Value = $arr
$ws.Range("M58").Value = -4678.2
$ws.Range("N58").Value = -7726.5

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 961.1667
$arr[0,1] = 1031.8889
$arr[0,2] = 749
$arr[0,3] = 1031.8889
$arr[0,4] = 749
$ws.Range("H107:L107").Value = $arr
$ws.Range("M107").Value = 888.1111000000001
$ws.Range("N107").Value = -4589

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 4381.25
$arr[0,1] = 1303.8334
$arr[0,2] = 13613.5
$arr[0,3] = 3911.5002
$arr[0,4] = 40840.5
$ws.Range("H132:L132").Value = $arr
$ws.Range("M132").Value = -1381.5002
$ws.Range("N132").Value = -45900.5

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 5965.3335
$arr[0,1] = 4881.2
$arr[0,2] = 7320.5
$arr[0,3] = 14643.6
$arr[0,4] = 21961.5
$ws.Range("H136:L136").Value = $arr
$ws.Range("M136").Value = -12093.6
$ws.Range("N136").Value = -27061.5

$ws = $wb.Worksheets.Item("CUL")
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1524.6774
$arr[0,1] = 791.1875
$arr[0,2] = 2307.0667
$arr[0,3] = 2373.5625
$arr[0,4] = 6921.2001
$ws.Range("H5:L5").Value = $arr
$ws.Range("M5").Value = -2261.5625
$ws.Range("N5").Value = -7145.2001

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 398.81818
$arr[0,1] = 354.1111
$arr[0,2] = 600
$arr[0,3] = 1062.3333
$arr[0,4] = 1800
$ws.Range("H34:L34").Value = $arr
$ws.Range("M34").Value = -978.3333
$ws.Range("N34").Value = -1968

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1909.1666
$arr[0,1] = 1909.1666
$arr[0,2] = 0
$arr[0,3] = 5727.4998
$arr[0,4] = 0
$ws.Range("H39:L39").Value = $arr
$ws.Range("M39").Value = -5433.4998
$ws.Range("N39").ClearContents()

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 999
$arr[0,1] = 0
$arr[0,2] = 999
$arr[0,3] = 0
$arr[0,4] = 2997
$ws.Range("H48:L48").Value = $arr
$ws.Range("N48").Value = -3497

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 4587.25
$arr[0,1] = 4587.25
$arr[0,2] = 0
$arr[0,3] = 13761.75
$arr[0,4] = 0
$ws.Range("H55:L55").Value = $arr
$ws.Range("M55").Value = -13584.75

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1397.1765
$arr[0,1] = 0
$arr[0,2] = 1397.1765
$arr[0,3] = 0
$arr[0,4] = 12574.5885
$ws.Range("H122:L122").Value = $arr
$ws.Range("N122").Value = -17474.5885

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1524.6774
$arr[0,1] = 791.1875
$arr[0,2] = 2307.0667
$arr[0,3] = 7120.6875
$arr[0,4] = 20763.6003
$ws.Range("H135:L135").Value = $arr
$ws.Range("M135").Value = -4585.6875
$ws.Range("N135").Value = -25833.6003

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 3118.1
$arr[0,1] = 2798.111
$arr[0,2] = 5998
$arr[0,3] = 8394.332999999999
$arr[0,4] = 17994
$ws.Range("H137:L137").Value = $arr
$ws.Range("M137").Value = -3294.332999999999
$ws.Range("N137").Value = -28194

$ws = $wb.Worksheets.Item("GSM")
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 4740
$arr[0,1] = 4263.4287
$arr[0,2] = 5073.6
$arr[0,3] = 4263.4287
$arr[0,4] = 5073.6
$ws.Range("H113:L113").Value = $arr
$ws.Range("M113").Value = -2093.4287
$ws.Range("N113").Value = -9413.6

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 50003276
$arr[0,1] = 66670250
$arr[0,2] = 2362.4
$arr[0,3] = 200010750
$arr[0,4] = 7087.200000000001
$ws.Range("H132:L132").Value = $arr
$ws.Range("M132").Value = -200008220
$ws.Range("N132").Value = -12147.2

$ws = $wb.Worksheets.Item("LTW")
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 311223.75
$arr[0,1] = 5017.885
$arr[0,2] = 1107359
$arr[0,3] = 5017.885
$arr[0,4] = 1107359
$ws.Range("H7:L7").Value = $arr
$ws.Range("M7").Value = -4905.885
$ws.Range("N7").Value = -1107583

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$ws.Range("H18:L18").Value = $arr
$ws.Range("M18").ClearContents()

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1188.4
$arr[0,1] = 1072.6
$arr[0,2] = 1420
$arr[0,3] = 1072.6
$arr[0,4] = 1420
$ws.Range("H22:L22").Value = $arr
$ws.Range("M22").Value = -777.5999999999999
$ws.Range("N22").Value = -2010

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1188.4
$arr[0,1] = 1072.6
$arr[0,2] = 1420
$arr[0,3] = 1072.6
$arr[0,4] = 1420
$ws.Range("H27:L27").Value = $arr
$ws.Range("M27").Value = -965.5999999999999
$ws.Range("N27").Value = -1634

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 3713.5557
$arr[0,1] = 2841.9
$arr[0,2] = 4803.125
$arr[0,3] = 2841.9
$arr[0,4] = 4803.125
$ws.Range("H46:L46").Value = $arr
$ws.Range("M46").Value = -2653.9
$ws.Range("N46").Value = -5179.125

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 5428.926
$arr[0,1] = 4851.3477
$arr[0,2] = 8750
$arr[0,3] = 14554.0431
$arr[0,4] = 26250
$ws.Range("H122:L122").Value = $arr
$ws.Range("M122").Value = -12104.0431
$ws.Range("N122").Value = -31150

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 311223.75
$arr[0,1] = 5017.885
$arr[0,2] = 1107359
$arr[0,3] = 15053.655
$arr[0,4] = 3322077
$ws.Range("H126:L126").Value = $arr
$ws.Range("M126").Value = -12583.655
$ws.Range("N126").Value = -3327017

$ws = $wb.Worksheets.Item("WVR")
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 16548.75
$arr[0,1] = 9000.5
$arr[0,2] = 18058.4
$arr[0,3] = 18001
$arr[0,4] = 36116.8
$ws.Range("H81:L81").Value = $arr
$ws.Range("M81").Value = -16940
$ws.Range("N81").Value = -38238.8

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 16548.75
$arr[0,1] = 9000.5
$arr[0,2] = 18058.4
$arr[0,3] = 90005
$arr[0,4] = 180584
$ws.Range("H84:L84").Value = $arr
$ws.Range("M84").Value = -84701
$ws.Range("N84").Value = -191192

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 930.13635
$arr[0,1] = 816.9375
$arr[0,2] = 1232
$arr[0,3] = 1633.875
$arr[0,4] = 2464
$ws.Range("H100:L100").Value = $arr
$ws.Range("M100").Value = -1092.875
$ws.Range("N100").Value = -3546

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 106999
$arr[0,1] = 0
$arr[0,2] = 106999
$arr[0,3] = 0
$arr[0,4] = 106999
$ws.Range("H109:L109").Value = $arr
$ws.Range("N109").Value = -109773

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1636
$arr[0,1] = 1458.4
$arr[0,2] = 1813.6
$arr[0,3] = 4375.200000000001
$arr[0,4] = 5440.799999999999
$ws.Range("H113:L113").Value = $arr
$ws.Range("M113").Value = -2205.200000000001
$ws.Range("N113").Value = -9780.799999999999

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 4298.6665
$arr[0,1] = 4298.6665
$arr[0,2] = 0
$arr[0,3] = 12895.9995
$arr[0,4] = 0
$ws.Range("H122:L122").Value = $arr
$ws.Range("M122").Value = -10445.9995

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 4508.3125
$arr[0,1] = 4212.2085
$arr[0,2] = 5396.625
$arr[0,3] = 12636.6255
$arr[0,4] = 16189.875
$ws.Range("H136:L136").Value = $arr
$ws.Range("M136").Value = -10086.6255
$ws.Range("N136").Value = -21289.875

